# Weekly update: insert the new week's two price rows (Melón - Calameño)
# for "Agrícola del Norte S.A. de Arica" at the top of the existing data
# block (old rows 33-37), pushing the previous rows down to 35-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 33; this shifts the
# existing rows 33-37 down to rows 35-39 (format/values carried along).
$ws.Rows("33:34").Insert()

# --- New row 33: Melón / Calameño / Primera -----------------------------
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C33").Value = 'Arica y Parinacota'
$ws.Range("D33").Value = 44559
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112027
$ws.Range("G33").Value = 'Melón'
$ws.Range("H33").Value = 'Calameño'
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = 8500
$ws.Range("N33").Value = '$/caja 16 unidades'
$ws.Range("O33").Value = 'Región de Arica y Parinacota'
$ws.Range("P33").Value = 531
$ws.Range("Q33").Value = 16
$ws.Range("R33").Value = 'Hortaliza'

# --- New row 34: Melón / Calameño / Super --------------------------------
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C34").Value = 'Arica y Parinacota'
$ws.Range("D34").Value = 44559
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100112027
$ws.Range("G34").Value = 'Melón'
$ws.Range("H34").Value = 'Calameño'
$ws.Range("I34").Value = 'Super'
$ws.Range("J34").Value = 70
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = 12500
$ws.Range("N34").Value = '$/caja 12 unidades'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 1042
$ws.Range("Q34").Value = 12
$ws.Range("R34").Value = 'Hortaliza'
